$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLineStyleNone
$xlNone = -4142
# xlEdgeBottom
$xlEdgeBottom = 9

# --- Row 3 (U=50 block + U=70 block: raw inputs in F/G and K/L, formulas elsewhere) ---
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 0.117
$ws.Range("E3").Formula = '=D3/$B$11'
$ws.Range("H3").Formula = '=F3*$B$10'
$ws.Range("I3").Formula = '=G3-$G$8'
$ws.Range("J3").Formula = '=I3/$B$12'
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 0.0971
$ws.Range("M3").Formula = '=K3*$B$10'
$ws.Range("N3").Formula = '=L3-$L$8'
$ws.Range("O3").Formula = '=N3/$B$13'

# --- Row 4 ---
$ws.Range("F4").Value = 3.5
$ws.Range("G4").Value = 0.1052
$ws.Range("E4").Formula = '=D4/$B$11'
$ws.Range("H4").Formula = '=F4*$B$10'
$ws.Range("I4").Formula = '=G4-$G$8'
$ws.Range("J4").Formula = '=I4/$B$12'
$ws.Range("K4").Value = 3.5
$ws.Range("L4").Value = 0.0872
$ws.Range("M4").Formula = '=K4*$B$10'
$ws.Range("N4").Formula = '=L4-$L$8'
$ws.Range("O4").Formula = '=N4/$B$13'

# --- Row 5 ---
$ws.Range("F5").Value = 2.7
$ws.Range("G5").Value = 0.0781
$ws.Range("E5").Formula = '=D5/$B$11'
$ws.Range("H5").Formula = '=F5*$B$10'
$ws.Range("I5").Formula = '=G5-$G$8'
$ws.Range("J5").Formula = '=I5/$B$12'
$ws.Range("K5").Value = 2.7
$ws.Range("L5").Value = 0.0659
$ws.Range("M5").Formula = '=K5*$B$10'
$ws.Range("N5").Formula = '=L5-$L$8'
$ws.Range("O5").Formula = '=N5/$B$13'

# --- Row 6 ---
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0.0521
$ws.Range("E6").Formula = '=D6/$B$11'
$ws.Range("H6").Formula = '=F6*$B$10'
$ws.Range("I6").Formula = '=G6-$G$8'
$ws.Range("J6").Formula = '=I6/$B$12'
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.041
$ws.Range("M6").Formula = '=K6*$B$10'
$ws.Range("N6").Formula = '=L6-$L$8'
$ws.Range("O6").Formula = '=N6/$B$13'

# --- Row 7 ---
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.0225
$ws.Range("E7").Formula = '=D7/$B$11'
$ws.Range("H7").Formula = '=F7*$B$10'
$ws.Range("I7").Formula = '=G7-$G$8'
$ws.Range("J7").Formula = '=I7/$B$12'
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.0189
$ws.Range("M7").Formula = '=K7*$B$10'
$ws.Range("N7").Formula = '=L7-$L$8'
$ws.Range("O7").Formula = '=N7/$B$13'

# --- Row 8 ---
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0.0095
$ws.Range("E8").Formula = '=D8/$B$11'
$ws.Range("H8").Formula = '=F8*$B$10'
$ws.Range("I8").Formula = '=G8-$G$8'
$ws.Range("J8").Formula = '=I8/$B$12'
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0.01
$ws.Range("M8").Formula = '=K8*$B$10'
$ws.Range("N8").Formula = '=L8-$L$8'
$ws.Range("O8").Formula = '=N8/$B$13'

# Row 8's newly-filled formula cells (E,H,I,J,M,N,O) pick up the plain
# (no-bottom-border) look of the rest of the table instead of the
# original outer-border treatment reserved for the table's last row.
foreach ($col in @("E", "H", "I", "J", "M", "N", "O")) {
    $ws.Range($col + "8").Borders.Item($xlEdgeBottom).LineStyle = $xlNone
}

# --- k-coefficients (rows 11-13) ---
$ws.Range("B11").Value = 40.8
$ws.Range("B12").Formula = '=(42.5+43.3)/2'
$ws.Range("B13").Formula = '=(44.1+44.9)/2'

# --- Selection moves to D10 ---
$ws.Range("D10").Select()
